$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lembar1")

# --- Fix test data (Login test-case sheet) ---
# TD001 valid credentials: admintitan/Testing123@ -> Admin/admin123
# Every cell that held the old "admintitan" / "Testing123@" text is updated
# (they were stored as shared strings reused across several rows).
$ws.Range("B2").Value = "Admin"
$ws.Range("B5").Value = "Admin"
$ws.Range("C5").Value = "Admin"
$ws.Range("B6").Value = "Admin"

$ws.Range("C2").Value = "admin123"
$ws.Range("C3").Value = "admin123"
$ws.Range("C4").Value = "admin123"

# Expected result text / validation object updated for the new login target
$ws.Range("F2").Value = "Success landing at Dashboard Page"
$ws.Range("G2").Value = "Dashboard Page"

# --- Capture new object for case logout: move the active selection ---
$ws.Range("F10").Select()

# --- Clear stray empty cells left over from earlier edits ---
$ws.Range("B4").ClearContents()
$ws.Range("C6").ClearContents()
